$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C1 header text from "job_prefix" to "deskripsi"
$ws.Range("C1").Value = "deskripsi"

# Remove the now-unused columns D:F (company_code, band, flag_mgr)
$ws.Range("D1:F2").ClearContents()
